# Apply data refresh for 2022-12-15 to CTA violent crime YTD workbook
# Updates numeric cell values across several worksheets as per the source diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("D2").Value = 92   # was 91
$ws.Range("H2").Value = 108   # was 107
$ws.Range("B3").Value = 77   # was 76
$ws.Range("I4").Value = 21   # was 22
$ws.Range("C6").Value = 489   # was 486
$ws.Range("D6").Value = 421   # was 419
$ws.Range("E6").Value = 486   # was 483
$ws.Range("F6").Value = 553   # was 551
$ws.Range("G6").Value = 439   # was 438
$ws.Range("I6").Value = 505   # was 502
$ws.Range("B7").Value = 517   # was 516
$ws.Range("C7").Value = 646   # was 643
$ws.Range("D7").Value = 660   # was 657
$ws.Range("E7").Value = 719   # was 716
$ws.Range("F7").Value = 801   # was 799
$ws.Range("G7").Value = 673   # was 672
$ws.Range("H7").Value = 727   # was 726
$ws.Range("I7").Value = 841   # was 839

# --- Sheet: Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("C6").Value = 35   # was 34
$ws.Range("C7").Value = 40   # was 39

# --- Sheet: Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("C4").Value = 8   # was 7
$ws.Range("C5").Value = 10   # was 9

# --- Sheet: Uptown ---
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("G5").Value = 14   # was 13
$ws.Range("G6").Value = 22   # was 21

# --- Sheet: Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("H2").Value = 11   # was 10
$ws.Range("B3").Value = 4   # was 3
$ws.Range("D6").Value = 24   # was 23
$ws.Range("F6").Value = 38   # was 37
$ws.Range("B7").Value = 36   # was 35
$ws.Range("D7").Value = 48   # was 47
$ws.Range("F7").Value = 58   # was 57
$ws.Range("H7").Value = 46   # was 45

# --- Sheet: By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I7").Value = 10   # was 9
$ws.Range("B28").Value = 36   # was 35
$ws.Range("D28").Value = 48   # was 47
$ws.Range("F28").Value = 58   # was 57
$ws.Range("H28").Value = 46   # was 45
$ws.Range("C36").Value = 40   # was 39
$ws.Range("C41").Value = 10   # was 9
$ws.Range("F47").Value = 18   # was 17
$ws.Range("D53").Value = 73   # was 72
$ws.Range("I53").Value = 124   # was 125
$ws.Range("E70").Value = 19   # was 18
$ws.Range("E74").Value = 7   # was 6
$ws.Range("I74").Value = 20   # was 19
$ws.Range("D76").Value = 15   # was 14
$ws.Range("C77").Value = 25   # was 24
$ws.Range("G86").Value = 22   # was 21
$ws.Range("E95").Value = 6   # was 5
$ws.Range("I95").Value = 6   # was 5
$ws.Range("B98").Value = 517   # was 516
$ws.Range("C98").Value = 646   # was 643
$ws.Range("D98").Value = 660   # was 657
$ws.Range("E98").Value = 719   # was 716
$ws.Range("F98").Value = 801   # was 799
$ws.Range("G98").Value = 673   # was 672
$ws.Range("H98").Value = 727   # was 726
$ws.Range("I98").Value = 841   # was 839

# --- Sheet: Loop ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I4").Value = 3   # was 4
$ws.Range("D6").Value = 44   # was 43
$ws.Range("D7").Value = 73   # was 72
$ws.Range("I7").Value = 124   # was 125

# --- Sheet: Rogers Park ---
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("D2").Value = 3   # was 2
$ws.Range("D7").Value = 15   # was 14

# --- Sheet: River North ---
$ws = $wb.Worksheets.Item("River North")
$ws.Range("E5").Value = 5   # was 4
$ws.Range("I5").Value = 13   # was 12
$ws.Range("E6").Value = 7   # was 6
$ws.Range("I6").Value = 20   # was 19

# --- Sheet: Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("C6").Value = 16   # was 15
$ws.Range("C7").Value = 25   # was 24

# --- Sheet: Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I5").Value = 5   # was 4
$ws.Range("I6").Value = 10   # was 9

# --- Sheet: Lake View ---
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("F5").Value = 11   # was 10
$ws.Range("F6").Value = 18   # was 17

# --- Sheet: Old Town ---
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("E4").Value = 17   # was 16
$ws.Range("E5").Value = 19   # was 18

# --- Sheet: Wicker Park ---
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("D4").Value = 3   # was 2
$ws.Range("G4").Value = 5   # was 4
$ws.Range("D5").Value = 6   # was 5
$ws.Range("G5").Value = 6   # was 5
